$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 52.813101
$ws.Range("H2").Value = 105.626202
$ws.Range("I2").Value = 0.2636577117692198
$ws.Range("J2").Value = 0.1954072982860194
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 5.5509805
$ws.Range("N2").Value = 11.101961
$ws.Range("O2").Value = 0.4587918645086687
$ws.Range("P2").Value = 0.3684119094562331
$ws.Range("Q2").Value = 293.1644937955305
$ws.Range("R2").Value = 1172.657975182122
$ws.Range("S2").Value = 0.1209640131746895
$ws.Range("T2").Value = 0.07199037588323609

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 52.813101
$ws.Range("H3").Value = 105.626202
$ws.Range("I3").Value = 0.2636577117692198
$ws.Range("J3").Value = 0.1954072982860194
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.327207
$ws.Range("N3").Value = 3.981621
$ws.Range("O3").Value = 0.1096944538210784
$ws.Range("P3").Value = 0.1321277020646205
$ws.Range("Q3").Value = 70.09391733890699
$ws.Range("R3").Value = 420.563504033442
$ws.Range("S3").Value = 0.02892178868823988
$ws.Range("T3").Value = 0.02581871728918759

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 52.813101
$ws.Range("H4").Value = 105.626202
$ws.Range("I4").Value = 0.2636577117692198
$ws.Range("J4").Value = 0.1954072982860194
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.250284
$ws.Range("N4").Value = 0.750852
$ws.Range("O4").Value = 0.02068612257180288
$ws.Range("P4").Value = 0.02491657276034671
$ws.Range("Q4").Value = 13.218274170684
$ws.Range("R4").Value = 79.309645024104
$ws.Range("S4").Value = 0.005454055742659156
$ws.Range("T4").Value = 0.004868880165646374

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 52.813101
$ws.Range("H5").Value = 105.626202
$ws.Range("I5").Value = 0.2636577117692198
$ws.Range("J5").Value = 0.1954072982860194
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.3570200000000001
$ws.Range("N5").Value = 1.07106
$ws.Range("O5").Value = 0.02950791692870925
$ws.Range("P5").Value = 0.03554248296694548
$ws.Range("Q5").Value = 18.85533331902
$ws.Range("R5").Value = 113.13199991412
$ws.Range("S5").Value = 0.007779989856499705
$ws.Range("T5").Value = 0.006945260570947678

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 52.813101
$ws.Range("H6").Value = 105.626202
$ws.Range("I6").Value = 0.2636577117692198
$ws.Range("J6").Value = 0.1954072982860194
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.001879333333334
$ws.Range("N6").Value = 12.005638
$ws.Range("O6").Value = 0.3307577248521605
$ws.Range("P6").Value = 0.3983998880756572
$ws.Range("Q6").Value = 211.3516574211461
$ws.Range("R6").Value = 1268.109944526876
$ws.Range("S6").Value = 0.08720682488451384
$ws.Range("T6").Value = 0.07785024576631668

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 52.813101
$ws.Range("H7").Value = 105.626202
$ws.Range("I7").Value = 0.2636577117692198
$ws.Range("J7").Value = 0.1954072982860194
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 0.611755
$ws.Range("N7").Value = 1.22351
$ws.Range("O7").Value = 0.05056191731758031
$ws.Range("P7").Value = 0.04060144467619692
$ws.Range("Q7").Value = 32.30867860225501
$ws.Range("R7").Value = 129.23471440902
$ws.Range("S7").Value = 0.01333103942261771
$ws.Range("T7").Value = 0.007933818610684923

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.139129333333334
$ws.Range("H8").Value = 12.417388
$ws.Range("I8").Value = 0.02066368662471691
$ws.Range("J8").Value = 0.02297202961864744
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 5.5509805
$ws.Range("N8").Value = 11.101961
$ws.Range("O8").Value = 0.4587918645086687
$ws.Range("P8").Value = 0.3684119094562331
$ws.Range("Q8").Value = 22.97622621631134
$ws.Range("R8").Value = 137.857357297868
$ws.Range("S8").Value = 0.009480331314176709
$ws.Range("T8").Value = 0.008463169295891044

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.139129333333334
$ws.Range("H9").Value = 12.417388
$ws.Range("I9").Value = 0.02066368662471691
$ws.Range("J9").Value = 0.02297202961864744
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.327207
$ws.Range("N9").Value = 3.981621
$ws.Range("O9").Value = 0.1096944538210784
$ws.Range("P9").Value = 0.1321277020646205
$ws.Range("Q9").Value = 5.493481425105333
$ws.Range("R9").Value = 49.441332825948
$ws.Range("S9").Value = 0.002266691818228243
$ws.Range("T9").Value = 0.003035241485272286

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.139129333333334
$ws.Range("H10").Value = 12.417388
$ws.Range("I10").Value = 0.02066368662471691
$ws.Range("J10").Value = 0.02297202961864744
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.250284
$ws.Range("N10").Value = 0.750852
$ws.Range("O10").Value = 0.02068612257180288
$ws.Range("P10").Value = 0.02491657276034671
$ws.Range("Q10").Value = 1.035957846064
$ws.Range("R10").Value = 9.323620614576001
$ws.Range("S10").Value = 0.0004274515543042177
$ws.Range("T10").Value = 0.0005723842474458686

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4.139129333333334
$ws.Range("H11").Value = 12.417388
$ws.Range("I11").Value = 0.02066368662471691
$ws.Range("J11").Value = 0.02297202961864744
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.3570200000000001
$ws.Range("N11").Value = 1.07106
$ws.Range("O11").Value = 0.02950791692870925
$ws.Range("P11").Value = 0.03554248296694548
$ws.Range("Q11").Value = 1.477751954586667
$ws.Range("R11").Value = 13.29976759128
$ws.Range("S11").Value = 0.0006097423483630268
$ws.Range("T11").Value = 0.0008164829714369437

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 4.139129333333334
$ws.Range("H12").Value = 12.417388
$ws.Range("I12").Value = 0.02066368662471691
$ws.Range("J12").Value = 0.02297202961864744
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 4.001879333333334
$ws.Range("N12").Value = 12.005638
$ws.Range("O12").Value = 0.3307577248521605
$ws.Range("P12").Value = 0.3983998880756572
$ws.Range("Q12").Value = 16.56429613706045
$ws.Range("R12").Value = 149.078665233544
$ws.Range("S12").Value = 0.006834673975049384
$ws.Range("T12").Value = 0.009152054028939823

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 4.139129333333334
$ws.Range("H13").Value = 12.417388
$ws.Range("I13").Value = 0.02066368662471691
$ws.Range("J13").Value = 0.02297202961864744
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 0.611755
$ws.Range("N13").Value = 1.22351
$ws.Range("O13").Value = 0.05056191731758031
$ws.Range("P13").Value = 0.04060144467619692
$ws.Range("Q13").Value = 2.532133065313334
$ws.Range("R13").Value = 15.19279839188
$ws.Range("S13").Value = 0.001044795614595327
$ws.Range("T13").Value = 0.0009326975896614709

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 45.51938266666667
$ws.Range("H14").Value = 136.558148
$ws.Range("I14").Value = 0.227245438116592
$ws.Range("J14").Value = 0.2526310541736829
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 5.5509805
$ws.Range("N14").Value = 11.101961
$ws.Range("O14").Value = 0.4587918645086687
$ws.Range("P14").Value = 0.3684119094562331
$ws.Range("Q14").Value = 252.6772055547047
$ws.Range("R14").Value = 1516.063233328228
$ws.Range("S14").Value = 0.1042583582546005
$ws.Range("T14").Value = 0.09307228905606758

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 45.51938266666667
$ws.Range("H15").Value = 136.558148
$ws.Range("I15").Value = 0.227245438116592
$ws.Range("J15").Value = 0.2526310541736829
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.327207
$ws.Range("N15").Value = 3.981621
$ws.Range("O15").Value = 0.1096944538210784
$ws.Range("P15").Value = 0.1321277020646205
$ws.Range("Q15").Value = 60.41364331087866
$ws.Range("R15").Value = 543.722789797908
$ws.Range("S15").Value = 0.02492756421753122
$ws.Range("T15").Value = 0.03337956065813138

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 45.51938266666667
$ws.Range("H16").Value = 136.558148
$ws.Range("I16").Value = 0.227245438116592
$ws.Range("J16").Value = 0.2526310541736829
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.250284
$ws.Range("N16").Value = 0.750852
$ws.Range("O16").Value = 0.02068612257180288
$ws.Range("P16").Value = 0.02491657276034671
$ws.Range("Q16").Value = 11.392773171344
$ws.Range("R16").Value = 102.534958542096
$ws.Range("S16").Value = 0.004700826986762868
$ws.Range("T16").Value = 0.006294700042841662

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 45.51938266666667
$ws.Range("H17").Value = 136.558148
$ws.Range("I17").Value = 0.227245438116592
$ws.Range("J17").Value = 0.2526310541736829
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.3570200000000001
$ws.Range("N17").Value = 1.07106
$ws.Range("O17").Value = 0.02950791692870925
$ws.Range("P17").Value = 0.03554248296694548
$ws.Range("Q17").Value = 16.25132999965334
$ws.Range("R17").Value = 146.26196999688
$ws.Range("S17").Value = 0.006705539510372534
$ws.Range("T17").Value = 0.008979134939889607

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 45.51938266666667
$ws.Range("H18").Value = 136.558148
$ws.Range("I18").Value = 0.227245438116592
$ws.Range("J18").Value = 0.2526310541736829
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 4.001879333333334
$ws.Range("N18").Value = 12.005638
$ws.Range("O18").Value = 0.3307577248521605
$ws.Range("P18").Value = 0.3983998880756572
$ws.Range("Q18").Value = 182.1630767598249
$ws.Range("R18").Value = 1639.467690838424
$ws.Range("S18").Value = 0.07516318409447639
$ws.Range("T18").Value = 0.1006481837072306

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 45.51938266666667
$ws.Range("H19").Value = 136.558148
$ws.Range("I19").Value = 0.227245438116592
$ws.Range("J19").Value = 0.2526310541736829
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 0.611755
$ws.Range("N19").Value = 1.22351
$ws.Range("O19").Value = 0.05056191731758031
$ws.Range("P19").Value = 0.04060144467619692
$ws.Range("Q19").Value = 27.84670994324667
$ws.Range("R19").Value = 167.08025965948
$ws.Range("S19").Value = 0.01148996505284844
$ws.Range("T19").Value = 0.01025718576952209

# Row 20
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 43.599874
$ws.Range("H20").Value = 130.799622
$ws.Range("I20").Value = 0.2176627161557187
$ws.Range("J20").Value = 0.2419778451548658
$ws.Range("K20").Value = 2
$ws.Range("M20").Value = 5.5509805
$ws.Range("N20").Value = 11.101961
$ws.Range("O20").Value = 0.4587918645086687
$ws.Range("P20").Value = 0.3684119094562331
$ws.Range("Q20").Value = 242.022050376457
$ws.Range("R20").Value = 1452.132302258742
$ws.Range("S20").Value = 0.09986188337910329
$ws.Range("T20").Value = 0.08914751997960882

# Row 21
$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 43.599874
$ws.Range("H21").Value = 130.799622
$ws.Range("I21").Value = 0.2176627161557187
$ws.Range("J21").Value = 0.2419778451548658
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 1.327207
$ws.Range("N21").Value = 3.981621
$ws.Range("O21").Value = 0.1096944538210784
$ws.Range("P21").Value = 0.1321277020646205
$ws.Range("Q21").Value = 57.86605797191799
$ws.Range("R21").Value = 520.7945217472619
$ws.Range("S21").Value = 0.02387639276591397
$ws.Range("T21").Value = 0.03197197663086098

# Row 22
$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 43.599874
$ws.Range("H22").Value = 130.799622
$ws.Range("I22").Value = 0.2176627161557187
$ws.Range("J22").Value = 0.2419778451548658
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 0.250284
$ws.Range("N22").Value = 0.750852
$ws.Range("O22").Value = 0.02068612257180288
$ws.Range("P22").Value = 0.02491657276034671
$ws.Range("Q22").Value = 10.912350864216
$ws.Range("R22").Value = 98.211157777944
$ws.Range("S22").Value = 0.004502597625708735
$ws.Range("T22").Value = 0.006029258585193124

# Row 23
$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 43.599874
$ws.Range("H23").Value = 130.799622
$ws.Range("I23").Value = 0.2176627161557187
$ws.Range("J23").Value = 0.2419778451548658
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 0.3570200000000001
$ws.Range("N23").Value = 1.07106
$ws.Range("O23").Value = 0.02950791692870925
$ws.Range("P23").Value = 0.03554248296694548
$ws.Range("Q23").Value = 15.56602701548
$ws.Range("R23").Value = 140.09424313932
$ws.Range("S23").Value = 0.006422773346800166
$ws.Range("T23").Value = 0.008600493439794989

# Row 24
$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 43.599874
$ws.Range("H24").Value = 130.799622
$ws.Range("I24").Value = 0.2176627161557187
$ws.Range("J24").Value = 0.2419778451548658
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 4.001879333333334
$ws.Range("N24").Value = 12.005638
$ws.Range("O24").Value = 0.3307577248521605
$ws.Range("P24").Value = 0.3983998880756572
$ws.Range("Q24").Value = 174.4814346965373
$ws.Range("R24").Value = 1570.332912268836
$ws.Range("S24").Value = 0.07199362478080711
$ws.Range("T24").Value = 0.09640394642648725

# Row 25
$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 43.599874
$ws.Range("H25").Value = 130.799622
$ws.Range("I25").Value = 0.2176627161557187
$ws.Range("J25").Value = 0.2419778451548658
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 0.611755
$ws.Range("N25").Value = 1.22351
$ws.Range("O25").Value = 0.05056191731758031
$ws.Range("P25").Value = 0.04060144467619692
$ws.Range("Q25").Value = 26.67244091887
$ws.Range("R25").Value = 160.03464551322
$ws.Range("S25").Value = 0.0110054442573854
$ws.Range("T25").Value = 0.009824650092920629

# Row 26
$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 46.66673533333333
$ws.Range("H26").Value = 140.000206
$ws.Range("I26").Value = 0.2329733422342775
$ws.Range("J26").Value = 0.2589988231702788
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 5.5509805
$ws.Range("N26").Value = 11.101961
$ws.Range("O26").Value = 0.4587918645086687
$ws.Range("P26").Value = 0.3684119094562331
$ws.Range("Q26").Value = 259.0461378339943
$ws.Range("R26").Value = 1554.276827003966
$ws.Range("S26").Value = 0.1068862740644804
$ws.Range("T26").Value = 0.09541825099107969

# Row 27
$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 46.66673533333333
$ws.Range("H27").Value = 140.000206
$ws.Range("I27").Value = 0.2329733422342775
$ws.Range("J27").Value = 0.2589988231702788
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 1.327207
$ws.Range("N27").Value = 3.981621
$ws.Range("O27").Value = 0.1096944538210784
$ws.Range("P27").Value = 0.1321277020646205
$ws.Range("Q27").Value = 61.93641780154731
$ws.Range("R27").Value = 557.4277602139259
$ws.Range("S27").Value = 0.02555588353126025
$ws.Range("T27").Value = 0.03422091934292993

# Row 28
$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 46.66673533333333
$ws.Range("H28").Value = 140.000206
$ws.Range("I28").Value = 0.2329733422342775
$ws.Range("J28").Value = 0.2589988231702788
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 0.250284
$ws.Range("N28").Value = 0.750852
$ws.Range("O28").Value = 0.02068612257180288
$ws.Range("P28").Value = 0.02491657276034671
$ws.Range("Q28").Value = 11.679937186168
$ws.Range("R28").Value = 105.119434675512
$ws.Range("S28").Value = 0.004819315113420846
$ws.Range("T28").Value = 0.006453363022366424

# Row 29
$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 46.66673533333333
$ws.Range("H29").Value = 140.000206
$ws.Range("I29").Value = 0.2329733422342775
$ws.Range("J29").Value = 0.2589988231702788
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 0.3570200000000001
$ws.Range("N29").Value = 1.07106
$ws.Range("O29").Value = 0.02950791692870925
$ws.Range("P29").Value = 0.03554248296694548
$ws.Range("Q29").Value = 16.66095784870667
$ws.Range("R29").Value = 149.94862063836
$ws.Range("S29").Value = 0.006874558029252811
$ws.Range("T29").Value = 0.00920546126098856

# Row 30
$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 46.66673533333333
$ws.Range("H30").Value = 140.000206
$ws.Range("I30").Value = 0.2329733422342775
$ws.Range("J30").Value = 0.2589988231702788
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 4.001879333333334
$ws.Range("N30").Value = 12.005638
$ws.Range("O30").Value = 0.3307577248521605
$ws.Range("P30").Value = 0.3983998880756572
$ws.Range("Q30").Value = 186.7546436846031
$ws.Range("R30").Value = 1680.791793161428
$ws.Range("S30").Value = 0.07705773262861339
$ws.Range("T30").Value = 0.103185102162766

# Row 31
$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 46.66673533333333
$ws.Range("H31").Value = 140.000206
$ws.Range("I31").Value = 0.2329733422342775
$ws.Range("J31").Value = 0.2589988231702788
$ws.Range("K31").Value = 2
$ws.Range("M31").Value = 0.611755
$ws.Range("N31").Value = 1.22351
$ws.Range("O31").Value = 0.05056191731758031
$ws.Range("P31").Value = 0.04060144467619692
$ws.Range("Q31").Value = 28.54860867384333
$ws.Range("R31").Value = 171.29165204306
$ws.Range("S31").Value = 0.01177957886724988
$ws.Range("T31").Value = 0.01051572639014818

# Row 32
$ws.Range("E32").Value = 2
$ws.Range("G32").Value = 7.571113
$ws.Range("H32").Value = 15.142226
$ws.Range("I32").Value = 0.03779710509947509
$ws.Range("J32").Value = 0.02801294959650559
$ws.Range("K32").Value = 2
$ws.Range("M32").Value = 5.5509805
$ws.Range("N32").Value = 11.101961
$ws.Range("O32").Value = 0.4587918645086687
$ws.Range("P32").Value = 0.3684119094562331
$ws.Range("Q32").Value = 42.0271006262965
$ws.Range("R32").Value = 168.108402505186
$ws.Range("S32").Value = 0.01734100432161828
$ws.Range("T32").Value = 0.01032030425034984

# Row 33
$ws.Range("E33").Value = 2
$ws.Range("G33").Value = 7.571113
$ws.Range("H33").Value = 15.142226
$ws.Range("I33").Value = 0.03779710509947509
$ws.Range("J33").Value = 0.02801294959650559
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 1.327207
$ws.Range("N33").Value = 3.981621
$ws.Range("O33").Value = 0.1096944538210784
$ws.Range("P33").Value = 0.1321277020646205
$ws.Range("Q33").Value = 10.048434171391
$ws.Range("R33").Value = 60.290605028346
$ws.Range("S33").Value = 0.004146132799904816
$ws.Range("T33").Value = 0.003701286658238321

# Row 34
$ws.Range("E34").Value = 2
$ws.Range("G34").Value = 7.571113
$ws.Range("H34").Value = 15.142226
$ws.Range("I34").Value = 0.03779710509947509
$ws.Range("J34").Value = 0.02801294959650559
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 0.250284
$ws.Range("N34").Value = 0.750852
$ws.Range("O34").Value = 0.02068612257180288
$ws.Range("P34").Value = 0.02491657276034671
$ws.Range("Q34").Value = 1.894928446092
$ws.Range("R34").Value = 11.369570676552
$ws.Range("S34").Value = 0.0007818755489470574
$ws.Range("T34").Value = 0.0006979866968532565

# Row 35
$ws.Range("E35").Value = 2
$ws.Range("G35").Value = 7.571113
$ws.Range("H35").Value = 15.142226
$ws.Range("I35").Value = 0.03779710509947509
$ws.Range("J35").Value = 0.02801294959650559
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 0.3570200000000001
$ws.Range("N35").Value = 1.07106
$ws.Range("O35").Value = 0.02950791692870925
$ws.Range("P35").Value = 0.03554248296694548
$ws.Range("Q35").Value = 2.703038763260001
$ws.Range("R35").Value = 16.21823257956
$ws.Range("S35").Value = 0.001115313837421004
$ws.Range("T35").Value = 0.0009956497838877021

# Row 36
$ws.Range("E36").Value = 2
$ws.Range("G36").Value = 7.571113
$ws.Range("H36").Value = 15.142226
$ws.Range("I36").Value = 0.03779710509947509
$ws.Range("J36").Value = 0.02801294959650559
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 4.001879333333334
$ws.Range("N36").Value = 12.005638
$ws.Range("O36").Value = 0.3307577248521605
$ws.Range("P36").Value = 0.3983998880756572
$ws.Range("Q36").Value = 30.29868064503134
$ws.Range("R36").Value = 181.792083870188
$ws.Range("S36").Value = 0.01250168448870037
$ws.Range("T36").Value = 0.01116035598391685

# Row 37
$ws.Range("E37").Value = 2
$ws.Range("G37").Value = 7.571113
$ws.Range("H37").Value = 15.142226
$ws.Range("I37").Value = 0.03779710509947509
$ws.Range("J37").Value = 0.02801294959650559
$ws.Range("K37").Value = 2
$ws.Range("M37").Value = 0.611755
$ws.Range("N37").Value = 1.22351
$ws.Range("O37").Value = 0.05056191731758031
$ws.Range("P37").Value = 0.04060144467619692
$ws.Range("Q37").Value = 4.631666233315
$ws.Range("R37").Value = 18.52666493326
$ws.Range("S37").Value = 0.001911094102883553
$ws.Range("T37").Value = 0.001137366223259614
